$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (C2): Mid Y changed ---
$ws.Cells.Item(3, 3).Value = -99.46

# --- Row 4: was J4 -> now J1, new coordinates ---
$ws.Cells.Item(4, 1).Value = "J1"
$ws.Cells.Item(4, 2).Value = 142.24
$ws.Cells.Item(4, 3).Value = -117.37
$ws.Cells.Item(4, 5).Value = 0

# --- Row 5: was Q1 -> now J3, new coordinates ---
$ws.Cells.Item(5, 1).Value = "J3"
$ws.Cells.Item(5, 2).Value = 130.66
$ws.Cells.Item(5, 3).Value = -72.090500000000006
$ws.Cells.Item(5, 5).Value = 90

# --- Row 6: was Q2 -> now Q1, new coordinates ---
$ws.Cells.Item(6, 1).Value = "Q1"
$ws.Cells.Item(6, 2).Value = 127.8475
$ws.Cells.Item(6, 3).Value = -93

# --- Row 7: was Q3 -> now Q2, new coordinates ---
$ws.Cells.Item(7, 1).Value = "Q2"
$ws.Cells.Item(7, 2).Value = 144.71
$ws.Cells.Item(7, 3).Value = -106.6
$ws.Cells.Item(7, 5).Value = 180

# --- Row 8: was R1 -> now Q3, new coordinates ---
$ws.Cells.Item(8, 1).Value = "Q3"
$ws.Cells.Item(8, 2).Value = 140.91999999999999
$ws.Cells.Item(8, 3).Value = -107.5
$ws.Cells.Item(8, 5).Value = -90

# --- Row 9: was R2 -> now R1, new coordinates ---
$ws.Cells.Item(9, 1).Value = "R1"
$ws.Cells.Item(9, 2).Value = 149.0275
$ws.Cells.Item(9, 3).Value = -106.7525

# --- Row 10: was R3 -> now R2, new coordinates ---
$ws.Cells.Item(10, 1).Value = "R2"
$ws.Cells.Item(10, 2).Value = 147.45750000000001
$ws.Cells.Item(10, 3).Value = -106.74250000000001
$ws.Cells.Item(10, 5).Value = 90

# --- Row 11: was R4 -> now R3, new coordinates ---
$ws.Cells.Item(11, 1).Value = "R3"
$ws.Cells.Item(11, 2).Value = 144.13
$ws.Cells.Item(11, 3).Value = -112.88
$ws.Cells.Item(11, 5).Value = -90

# --- Row 12: was R5 -> now R4, new coordinates ---
$ws.Cells.Item(12, 1).Value = "R4"
$ws.Cells.Item(12, 2).Value = 142.31
$ws.Cells.Item(12, 3).Value = -112.9
$ws.Cells.Item(12, 5).Value = -90

# --- Row 13: was U1 -> becomes the new R5 row (brand-new component entry) ---
$ws.Cells.Item(13, 1).Value = "R5"
$ws.Cells.Item(13, 2).Value = 140.93
$ws.Cells.Item(13, 3).Value = -104.22
$ws.Cells.Item(13, 5).Value = 180

# --- Row 14: brand-new row carrying what used to be row 13's U1 data ---
$ws.Cells.Item(13, 1).Copy()
$ws.Cells.Item(14, 1).PasteSpecial(-4122)
$ws.Cells.Item(13, 2).Copy()
$ws.Cells.Item(14, 2).PasteSpecial(-4122)
$ws.Cells.Item(13, 3).Copy()
$ws.Cells.Item(14, 3).PasteSpecial(-4122)
$ws.Cells.Item(13, 4).Copy()
$ws.Cells.Item(14, 4).PasteSpecial(-4122)
$ws.Cells.Item(13, 5).Copy()
$ws.Cells.Item(14, 5).PasteSpecial(-4122)

$ws.Cells.Item(14, 1).Value = "U1"
$ws.Cells.Item(14, 2).Value = 140.85
$ws.Cells.Item(14, 3).Value = -89.84
$ws.Cells.Item(14, 4).Value = "Top"
$ws.Cells.Item(14, 5).Value = 0

# --- Row 15: new blank separator row (was row 14, now gains styled C/E cells) ---
$ws.Cells.Item(13, 1).Copy()
$ws.Cells.Item(15, 1).PasteSpecial(-4122)
$ws.Cells.Item(13, 2).Copy()
$ws.Cells.Item(15, 2).PasteSpecial(-4122)
$ws.Cells.Item(13, 3).Copy()
$ws.Cells.Item(15, 3).PasteSpecial(-4122)
$ws.Cells.Item(13, 5).Copy()
$ws.Cells.Item(15, 5).PasteSpecial(-4122)

# --- Row 19 -> row 18: trailing footnote cell moves up one row ---
$ws.Cells.Item(19, 6).Copy()
$ws.Cells.Item(18, 6).PasteSpecial(-4122)
$ws.Cells.Item(18, 6).Value = "``"
$ws.Cells.Item(19, 6).Value = ""

# --- Sheet view / selection ---
$ws.Range("A5:F5").Select()

# --- Column widths (minor autosize drift) ---
$ws.Columns.Item(1).ColumnWidth = 13.5546875
$ws.Columns.Item(2).ColumnWidth = 16.33203125
$ws.Columns.Item(3).ColumnWidth = 14.44140625
$ws.Columns.Item(4).ColumnWidth = 12.33203125
$ws.Columns.Item(5).ColumnWidth = 14.109375
